# POD_Excluded.xlsx update: append 5 new dispatch rows (57-61) to Table1 and
# extend the trailing blank I:J filler block by 5 rows (now ending at 100),
# matching a fresh export/paste of newer CM/PCM/TT tickets into the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- 1. Clone the formatting of the last existing table row onto the five
#        new rows so the new records look identical to the existing ones
#        (borders, fonts, fill, number formats) before any values go in.
$ws.Range("A56:J56").Copy()
$ws.Range("A57:J61").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- 2. Fill in the new values. Column order mirrors how the source data
#        was pasted in: Task/Order/TT ids first, then Region, then
#        Site ID/Title pairs, then Created At timestamps.
$ws.Range("B57").Value = "CM20251004000621"
$ws.Range("C57").Value = "PCM-20251004-00001009"
$ws.Range("D57").Value = "TT-20251004-00809"
$ws.Range("B58").Value = "CM20251005000096"
$ws.Range("C58").Value = "PCM-20251005-00000723"
$ws.Range("D58").Value = "TT-20251004-01712"
$ws.Range("B59").Value = "CM20251005000138"
$ws.Range("C59").Value = "PCM-20251005-00000797"
$ws.Range("D59").Value = "TT-20251005-00852"
$ws.Range("B60").Value = "CM20251005000645"
$ws.Range("C60").Value = "PCM-20251005-00001592"
$ws.Range("D60").Value = "TT-20251005-01568"
$ws.Range("B61").Value = "CM20251007000454"
$ws.Range("C61").Value = "PCM-20251007-00000864"
$ws.Range("D61").Value = "TT-20251007-00353"

$ws.Range("F57").Value = "Region_6"
$ws.Range("F58").Value = "Region_6"
$ws.Range("F59").Value = "Region_6"
$ws.Range("F60").Value = "Region_6"
$ws.Range("F61").Value = "Region_6"

$ws.Range("H57").Value = "WDS0920"
$ws.Range("J57").Value = "1(2G/5G) sites down under WDS0920-P2-USF/WDS0920-P2-USF @ Generator_SG"
$ws.Range("H58").Value = "RYAB0648"
$ws.Range("J58").Value = "BASE STATION SERVICE PROBLEM (Rx signal level failure) at RYAB0648-GSM-UMTS-P1-HUB-USF (Ref:TT-20250927-00484) @ Sceco_STB"
$ws.Range("H59").Value = "FIF0045"
$ws.Range("J59").Value = "BASE STATION SERVICE PROBLEM (TX out of order)  at FIF1M7045-P2-CA-L800 @ Sceco"
$ws.Range("H60").Value = "SMT0300"
$ws.Range("J60").Value = "CELL SERVICE PROBLEM (Cell disabled due to unknown problem; logs collected)  at SMT1MT0300-P3 @ Sceco"
$ws.Range("H61").Value = "FIF0045"
$ws.Range("J61").Value = "BASE STATION SERVICE PROBLEM (TX out of order)  at FIF1M7045-P2-CA-L800 @ Sceco"

$ws.Range("I57").Value = "2025-10-04 12:33:20"
$ws.Range("I58").Value = "2025-10-05 01:48:08"
$ws.Range("I59").Value = "2025-10-05 02:57:45"
$ws.Range("I60").Value = "2025-10-05 10:26:52"
$ws.Range("I61").Value = "2025-10-07 09:20:50"

# --- 3. Grow the Excel Table (ListObject) so the new rows participate in
#        the table (autofilter, banding, structured refs, etc.)
$lo.Resize($ws.Range("A1:J61"))

# --- 4. Re-create the trailing blank I:J "buffer" rows below the table so
#        the used range keeps the same 39-row cushion it had before
#        (was rows 57:95, now pushed down to rows 62:100).
$ws.Range("I96:J100").Value = "x"
$ws.Range("I96:J100").ClearContents()
$ws.Range("I96:J100").Style = "Normal"

# --- 5. Restore the active selection to where the user clicked next (C2).
$ws.Range("C2").Select()
